# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF (header "Date") holds the literal date string "4-24-2012-13"
# for every data row; it needs to become "2013-04-24" -- kept as literal
# text (not auto-converted to an Excel date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 2) { $lastRow = 31 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # Column BF
    if ($cell.Value2 -eq "4-24-2012-13") {
        # Put the replacement text in via a formula (so it's treated as a
        # string literal), then flatten the formula down to a static value
        # with Copy/PasteSpecial so Excel's "looks like a date" auto-detect
        # (which would otherwise turn "2013-04-24" into a date serial) never
        # gets a chance to run against the raw Value setter.
        $cell.Formula = '="2013-04-24"'
        $cell.Copy() | Out-Null
        $cell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
    }
}

$excel.CutCopyMode = 0
